$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Clear()

$ws.Cells.Item(1,1).Value = "MOOD     "
$ws.Cells.Item(1,2).Value = "EMOTION  "
$ws.Cells.Item(1,3).Value = "INTENSITY"
$ws.Cells.Item(1,4).Value = "   EVENT    "
$ws.Cells.Item(1,5).Value = " APPLIED STRATEGY    "
$ws.Cells.Item(1,6).Value = " PERSONALITY TRAITS "
$ws.Cells.Item(1,7).Value = " STRATEGIES RELATED "
$ws.Cells.Item(1,8).Value = " DOMINANT PERSONALITY "

$ws.Cells.Item(2,1).Value = -1.4827710390090942
$ws.Cells.Item(2,2).Value = "Distress"
$ws.Cells.Item(2,3).Value = 4.774208068847656
$ws.Cells.Item(2,4).Value = "Talk"
$ws.Cells.Item(2,5).Value = "None"

$ws.Cells.Item(3,1).Value = -0.974892795085907
$ws.Cells.Item(3,2).Value = "Love"
$ws.Cells.Item(3,3).Value = 1.5804238319396973
$ws.Cells.Item(3,4).Value = "Hello"
$ws.Cells.Item(3,5).Value = "None"

$ws.Cells.Item(4,1).Value = 0
$ws.Cells.Item(4,2).Value = "Love"
$ws.Cells.Item(4,3).Value = 2.3942959308624268
$ws.Cells.Item(4,4).Value = "Conversation"
$ws.Cells.Item(4,5).Value = "None"

$ws.Cells.Item(5,1).Value = 1.4531155824661255
$ws.Cells.Item(5,2).Value = "Love"
$ws.Cells.Item(5,3).Value = 4.6787238121032715
$ws.Cells.Item(5,4).Value = "Hug"
$ws.Cells.Item(5,5).Value = "None"

$ws.Cells.Item(6,1).Value = 1.419926404953003
$ws.Cells.Item(6,2).Value = "--"
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = "Not-Discussion"
$ws.Cells.Item(6,5).Value = "Situation Selection"

$ws.Cells.Item(7,1).Value = 2.123051166534424
$ws.Cells.Item(7,2).Value = "Joy"
$ws.Cells.Item(7,3).Value = 2.316424608230591
$ws.Cells.Item(7,4).Value = "Congrat"
$ws.Cells.Item(7,5).Value = "None"

$ws.Cells.Item(8,1).Value = 0.8047746419906616
$ws.Cells.Item(8,2).Value = "Distress"
$ws.Cells.Item(8,3).Value = 4.166054725646973
$ws.Cells.Item(8,4).Value = "Bye"
$ws.Cells.Item(8,5).Value = "None"

$ws.Cells.Item(9,1).Value = -0.793574869632721
$ws.Cells.Item(9,2).Value = "Hate"
$ws.Cells.Item(9,3).Value = 5.116583347320557
$ws.Cells.Item(9,4).Value = "Fired"
$ws.Cells.Item(9,5).Value = "None"

$ws.Cells.Item(10,1).Value = -2.7233526706695557
$ws.Cells.Item(10,2).Value = "Hate"
$ws.Cells.Item(10,3).Value = 6.242823123931885
$ws.Cells.Item(10,4).Value = "Crash"
$ws.Cells.Item(10,5).Value = "None"

$ws.Cells.Item(11,1).Value = -0.5619249939918518
$ws.Cells.Item(11,2).Value = "Joy"
$ws.Cells.Item(11,3).Value = 6.858622074127197
$ws.Cells.Item(11,4).Value = "Profits"
$ws.Cells.Item(11,5).Value = "None"

$ws.Cells.Item(12,6).Value = "Low Conscientiousness"

$ws.Cells.Item(13,6).Value = "Low Extraversion"

$ws.Cells.Item(14,6).Value = "High Neuroticism"

$ws.Cells.Item(15,6).Value = "Low Agreeableness"

$ws.Cells.Item(16,6).Value = "Low Openness"

$ws.Cells.Item(17,7).Value = "[Situation Selection, Strongly]"

$ws.Cells.Item(18,7).Value = "[Situation Modification, Weakly]"

$ws.Cells.Item(19,7).Value = "[Attention Deployment, Weakly]"

$ws.Cells.Item(20,7).Value = "[Cognitive Change, Weakly]"

$ws.Cells.Item(21,7).Value = "[Response Modulation, Lightly]"

$ws.Cells.Item(22,8).Value = "Neuroticism"
